$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: merge "Testowanie " + "jednostkowe" into one run "Testowanie jednostkowe" ---
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$mergedTitle = $titleRange.Characters(1, 22)
$mergedTitle.Text = "Testowanie jednostkowe"

# --- Subtitle shape: reposition and rebuild the text ---
$subShape = $s.Shapes.Item(2)

# Resize / reposition the subtitle placeholder (EMU -> points, 12700 EMU per point)
$subShape.Left = 99.1836220472441
$subShape.Top = 306.0
$subShape.Width = 527.3027041653544
$subShape.Height = 138.0

$subRange = $subShape.TextFrame.TextRange

# Prepend a new paragraph for the GitHub link, before the existing "Robert Pająk" paragraph.
$null = $subRange.InsertBefore("https://github.com/Pellared/Examples" + [char]13)

# Prepend a new paragraph for the name, before the link paragraph.
$afterLink = $subShape.TextFrame.TextRange
$null = $afterLink.InsertBefore("Robert Pajak" + [char]13)

# Remove the original "Robert Pajak" text that is now trailing (3rd paragraph),
# leaving just its paragraph mark / endParaRPr behind as an empty final paragraph.
$whole = $subShape.TextFrame.TextRange
$trailStart = $whole.Length - 12 + 1
$trailing = $whole.Characters($trailStart, 12)
$trailing.Text = ""

# Split paragraph 1 into "Robert " + "Pająk" runs.
$whole = $subShape.TextFrame.TextRange
$nameFirst = $whole.Characters(1, 7)
$nameFirst.Text = "Robert "
$nameLast = $whole.Characters(8, 5)
$nameLast.Text = "Pająk"

# Split paragraph 2 into "https://" + "github.com/Pellared/Examples" runs
# (position 13 is the paragraph mark ending paragraph 1, so paragraph 2 starts at 14).
$whole = $subShape.TextFrame.TextRange
$linkPrefix = $whole.Characters(14, 8)
$linkPrefix.Text = "https://"

# Apply the hyperlink across the whole URL text (both runs).
$whole = $subShape.TextFrame.TextRange
$urlRange = $whole.Characters(14, 37)
$action = $urlRange.ActionSettings(1)
$action.Hyperlink.Address = "https://github.com/Pellared/Examples"
